$d = $word.ActiveDocument

function Set-ParaXml($para, [string]$bodyXml) {
    $pkg = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body>' + $bodyXml + '</w:body>' +
        '</w:document>' +
        '</pkg:xmlData>' +
        '</pkg:part>' +
        '</pkg:package>'
    $para.Range.InsertXML($pkg)
}

# ---------------------------------------------------------------------
# Hunk 1: insert 5 new paragraphs right after paragraph 1 ("LAB EXERCISE -1")
# ---------------------------------------------------------------------

$p1 = $d.Paragraphs.Item(1)
$p1.Range.InsertParagraphAfter()
$d = $word.ActiveDocument
$p2 = $d.Paragraphs.Item(2)
$p2.Range.InsertParagraphAfter()
$d = $word.ActiveDocument
$p3 = $d.Paragraphs.Item(3)
$p3.Range.InsertParagraphAfter()
$d = $word.ActiveDocument
$p4 = $d.Paragraphs.Item(4)
$p4.Range.InsertParagraphAfter()
$d = $word.ActiveDocument
$p5 = $d.Paragraphs.Item(5)
$p5.Range.InsertParagraphAfter()
$d = $word.ActiveDocument

# paragraph 2: blank line, big bold title-style run properties, no jc
$p2 = $d.Paragraphs.Item(2)
Set-ParaXml $p2 @'
<w:p>
<w:pPr>
<w:pStyle w:val="Standard"/>
<w:rPr>
<w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
<w:b/>
<w:bCs/>
<w:color w:val="333333"/>
<w:sz w:val="40"/>
<w:szCs w:val="40"/>
<w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
</w:rPr>
</w:pPr>
</w:p>
'@

# paragraph 3: right aligned "N. Jyothi Kumar"
$d = $word.ActiveDocument
$p3 = $d.Paragraphs.Item(3)
Set-ParaXml $p3 @'
<w:p>
<w:pPr>
<w:pStyle w:val="Standard"/>
<w:jc w:val="right"/>
<w:rPr>
<w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
<w:b/>
<w:bCs/>
<w:color w:val="333333"/>
<w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
</w:rPr>
</w:pPr>
<w:r>
<w:rPr>
<w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
<w:b/>
<w:bCs/>
<w:color w:val="333333"/>
<w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
</w:rPr>
<w:t>N. Jyothi Kumar</w:t>
</w:r>
</w:p>
'@

# paragraph 4: right aligned "Ch.en.u4cse22139"
$d = $word.ActiveDocument
$p4 = $d.Paragraphs.Item(4)
Set-ParaXml $p4 @'
<w:p>
<w:pPr>
<w:pStyle w:val="Standard"/>
<w:jc w:val="right"/>
<w:rPr>
<w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
<w:b/>
<w:bCs/>
<w:color w:val="333333"/>
<w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
</w:rPr>
</w:pPr>
<w:r>
<w:rPr>
<w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
<w:b/>
<w:bCs/>
<w:color w:val="333333"/>
<w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
</w:rPr>
<w:t>Ch.</w:t>
</w:r>
<w:proofErr w:type="gramStart"/>
<w:r>
<w:rPr>
<w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
<w:b/>
<w:bCs/>
<w:color w:val="333333"/>
<w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
</w:rPr>
<w:t>en.u</w:t>
</w:r>
<w:proofErr w:type="gramEnd"/>
<w:r>
<w:rPr>
<w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
<w:b/>
<w:bCs/>
<w:color w:val="333333"/>
<w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
</w:rPr>
<w:t>4cse22139</w:t>
</w:r>
</w:p>
'@

# paragraph 5: right aligned "Cse - b"
$d = $word.ActiveDocument
$p5 = $d.Paragraphs.Item(5)
Set-ParaXml $p5 @'
<w:p>
<w:pPr>
<w:pStyle w:val="Standard"/>
<w:jc w:val="right"/>
<w:rPr>
<w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
<w:b/>
<w:bCs/>
<w:color w:val="333333"/>
<w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
</w:rPr>
</w:pPr>
<w:proofErr w:type="spellStart"/>
<w:r>
<w:rPr>
<w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
<w:b/>
<w:bCs/>
<w:color w:val="333333"/>
<w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
</w:rPr>
<w:t>Cse</w:t>
</w:r>
<w:proofErr w:type="spellEnd"/>
<w:r>
<w:rPr>
<w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
<w:b/>
<w:bCs/>
<w:color w:val="333333"/>
<w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
</w:rPr>
<w:t xml:space="preserve"> - b</w:t>
</w:r>
</w:p>
'@

# paragraph 6: blank right-aligned paragraph
$d = $word.ActiveDocument
$p6 = $d.Paragraphs.Item(6)
Set-ParaXml $p6 @'
<w:p>
<w:pPr>
<w:pStyle w:val="Standard"/>
<w:jc w:val="right"/>
<w:rPr>
<w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
<w:b/>
<w:bCs/>
<w:color w:val="333333"/>
<w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
</w:rPr>
</w:pPr>
</w:p>
'@

# ---------------------------------------------------------------------
# Hunk 2: insert 2 new paragraphs after the (now shifted) paragraph that
# used to be #26 (empty paragraph right before the final empty paragraph,
# just after the last "Output:" image). Originally index 26, now +5 = 31.
# ---------------------------------------------------------------------

$d = $word.ActiveDocument
$target = $d.Paragraphs.Item(31)
$target.Range.InsertParagraphAfter()
$d = $word.ActiveDocument
$next = $d.Paragraphs.Item(32)
$next.Range.InsertParagraphAfter()
$d = $word.ActiveDocument

# paragraph 32: blank paragraph
$p32 = $d.Paragraphs.Item(32)
Set-ParaXml $p32 @'
<w:p>
<w:pPr>
<w:pStyle w:val="Standard"/>
<w:rPr>
<w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
</w:rPr>
</w:pPr>
</w:p>
'@

# paragraph 33: "Result: ..." paragraph
$d = $word.ActiveDocument
$p33 = $d.Paragraphs.Item(33)
Set-ParaXml $p33 @'
<w:p>
<w:pPr>
<w:pStyle w:val="Standard"/>
<w:rPr>
<w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
</w:rPr>
</w:pPr>
<w:r>
<w:rPr>
<w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
</w:rPr>
<w:lastRenderedPageBreak/>
<w:t>Result:</w:t>
</w:r>
<w:r>
<w:t xml:space="preserve"> </w:t>
</w:r>
<w:r>
<w:rPr>
<w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
</w:rPr>
<w:t>The program successfully counts characters, words, spaces, lines, and distinguishes between positive/negative integers and fractions from the input file.</w:t>
</w:r>
<w:r>
<w:rPr>
<w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
</w:rPr>
<w:br/>
<w:t>All tasks executed correctly with accurate outputs based on the file content.</w:t>
</w:r>
</w:p>
'@

Write-Host "Edit complete. Paragraph count:" $d.Paragraphs.Count
